$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim leading whitespace from "name" column (column C) values that had an
# accidental leading space in the shared strings table for kan/hin/tam rows.
$ws.Range("C23").Value = "ಸೂರ್ಯ"
$ws.Range("C26").Value = "ಬುಧ"
$ws.Range("C32").Value = "मंगल"
$ws.Range("C34").Value = "इकट्ठा करना"
$ws.Range("C37").Value = "சூரியன்"

# Update the saved selection/active cell for the sheet view.
$ws.Range("B27").Select()
